# Daily attendance processing - 2026-01-31 09:43:37
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$col = $ws.Columns.Item(7)  # Column G - "Recorded By"

$firstAddress = $null
$found = $col.Find($oldValue)
while ($found -ne $null) {
    if ($firstAddress -eq $null) {
        $firstAddress = $found.Address()
    } elseif ($found.Address() -eq $firstAddress) {
        break
    }
    $found.Value = $newValue
    $found = $col.Find($oldValue)
}
